$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" everywhere it occurs ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"

# --- 2. Narrow the "status" columns (previously auto-sized wider for "Ready for handoff") ---
# Target stored width ~13.41 chars; engine quantizes ColumnWidth to 1/6-character
# pixel steps, so 12.5 is the closest input that lands on the nearest representable width.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
